$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ViewSwap")
Write-Output "before insert"
$ws.Range("B8").Insert(-4121)
Write-Output "after insert"
